$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-16
# from serial date 45185 (2023-09-16) to 45204 (2023-10-05)
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
